$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a temporary blank column before column A. This shifts the original
# column A (value "FICHA", custom width 10.7109375) into column B, carrying
# its exact custom width along with it (avoids lossy character-width
# rounding that happens when setting ColumnWidth directly).
$ws.Columns.Item(1).Insert()

# The former column B ("NOMBRE") is now in column C. Copy its value into the
# now-empty column A.
$ws.Range("A1").Value2 = $ws.Range("C1").Value2

# Remove the now-duplicate "NOMBRE" cell (old column B, currently column C),
# which shifts the remaining columns (D, E, ...) back down to (C, D, ...),
# restoring the original column count/layout.
$ws.Columns.Item(3).Delete()
